$d = $word.ActiveDocument

function Replace-WholeText($oldText, $newText, $label) {
    $ok = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $label"
    }
}

# --- 1. Title (Heading1 paragraph) ---
Replace-WholeText 'Review 173: [Short] Teaching Language Models to Self-Improve through Interactive Demonstrations' 'Review 172: [Short] Localizing and Editing Knowledge in Text-to-Image Generative Models' 'title'

# --- 2. "Paper: ..." bold link line ---
Replace-WholeText 'Paper: https://arxiv.org/abs/1909.04157v1' 'Paper: https://arxiv.org/abs/2405.01008v2' 'paper link'

# --- 3. huggingface.co paper link ---
Replace-WholeText 'https://huggingface.co/papers/2310.13522' 'https://huggingface.co/papers/2310.13730' 'huggingface link'

# --- 4. Paragraph 7: diffusion-models intro ---
Replace-WholeText 'ממודלי דיפוזיה שסקרנו אתמול עוברים לאייטם פופולרי אפילו אפילו מהם כלומר מודלי שפה ענקיים (LLMs). המאמר שנסקור היום מציע שיטה לאימון LLMs קטנים יחסית (מיליארדי פרמטרים בודדים) לפתרון בעיות מתמטיות מורכבות (נניח כאלו שמכילות הרבה פעולות).' 'מודלי דיפוזיה ממשיכים לשלוט ב AI גנרטיבי כבר זמן מה ואחד נושאי המחקר החמים ביותר בנושא הזה הוא עריכת תמונות המגונרטות עם מודלים אלו. לאחרונה יצאו לא מעט שיטות שמצליחות למשל להוריד אובייקט מתמונה, להחליף אותו לאובייקט אחר או לשנות את סגנון התמונה. המאמר שנסקור היום מציע שיטה לעריכת תמונות המגונרטות עם מודלי דיפוזיה בצורה מאוד אלגנטית המתבסס על ההבנה של מה שקורה בתוך מודל הדיפוזיה (שזה אנקודר של טקסט ומודל המסיר רעש מתמונה UNet בכל איטרציה). ' 'paragraph 7 (diffusion intro)'

# --- 5. Paragraph 9: "few-shot / small model" -> "first stage" paragraph ---
Replace-WholeText ' המאמר מציין שמודלי קטנים יחסית מתקשים לפתור בעיות בעזרת reasoning אם מפעילים אותו (המודל) בצורה של few-shot, כלומר מספקים לו כמה דוגמאות עם פתרון מלא.  בגדול המאמר מציע לאמן (מכייל) מודל שפה קטן L על הטעויות שלו. עבור בעיה נתונה מפעילים מודל L כדי ליצור שרשרת צעדי חשיבה לפתרון בעיה זו. לאחרי מכן מפעילים מודל יותר חזק (נגיד codex) לפתרון בעיה זו ומשווים את שרשרת החשיבה של שניהם. ' 'כלומר בשלב הראשון המאמר מנסה להבין איזה חלק(שכבה) במודל להסרת הרעש אחראי על יצירה של כל אובייקט בתמונה, איזו שכבה אחראית על הסגנון, ואיזו מהשכבות אחראית על צבע. איך עושים זאת? קודם כל מוסיפים את הרעש לטוקן האחרון של האובייקט/סגנון/צבע בתיאור הטקסטואלי. למה אותו דווקא?' 'paragraph 9 (first stage)'

# --- 6. Paragraph 11: "first place they differ" -> "paper checked and found" paragraph ---
Replace-WholeText 'במקום הראשון שהם שונים מחליפים את המשוב של המודל החלש בזה של המודל החזק. לאחר מכן מפעילים מודל חזק שוב פעם כדי לתקן את שרשרת החשיבה של המודל החלש מהמקום הזה.  לאחר מכן מחלקים את הדאטהסט הזה (יש בו פתרונות זהב ground-truth, פתרונות נכונים של המודל החלש, והפתרונות המתוקנים על ידי המודל החזק). אז מחלקים את הפתרונות האלו לפי התוצאה הסופית (נכונה או לא נכונה).' 'המאמר בדק ומצא (על ידי השימוש ב Clip-Score המודד את איכות התמונה המגונרטת והתאמתה לתיאור) שזה מה שמשפיע על הישות שרוצים לערוך (למשל מעלים אובייקט). אז איך עושים עריכה? מכיוון שהשכבה הראשונה אחרי שכבת האמבדינג באנקודר היא קריטית אז מאמנים רק אותה (את חלקה). מכיוון שיש לנו טרנספורמרים כאן אז השכבה מוגדרת על ידי 4 מטריצות: W_q, W_k, W_v ו- W_out. ' 'paragraph 11 (clip-score)'

# --- 7. Insert two brand-new paragraphs right after paragraph 11 (the one just edited above) ---
$p11 = $d.Paragraphs(11)
$rng = $p11.Range
$rng.InsertParagraphAfter()
$d.Paragraphs(12).Range.Text = ' '
$d.Paragraphs(12).Range.InsertParagraphAfter()
$d.Paragraphs(13).Range.Text = 'שלוש המטריצות הראשונות הן מטריצות ממנגנון ה-attention ומשאירים אותן כמו שהן ומאמנים רק את W_out (לצורך עריכה) תוך כדי שימוש בשיכונים (embeddings) של האובייקט (או סגנון) הישן והחדש c_k ו-c_v בהתאמה. פונקצייה שמאפטמים אותה כדי למצוא את W_out מרמזת על כך שהמטרה(לא לגמרי הפנמתי מה הרציונל כאן) היא למצוא W_out חדשה כך שפלט של השכבה הראשונה ״החדשה״ עבור c_k (הישן) תהיה כמה שיותר קרובה לפלט של השכבה המקורית עם c_v (החדש) עם רגולריזציה קטנה. והכי כיף שניתן לפתור בעיה זו בצורה סגורה ואין צורך באימון שזה מגניב. לבסוף הם עשו עוד דבר נחמד.'

# --- 8. Paragraph that was at index 13 (now 15 after the two insertions): "correct solutions divided" text ---
Replace-WholeText 'את הפתרונות הנכונים מחלקים לשלישיות של (תוצאה של שלב i, המשוב והתוצאה של השלב החדש). אלו שמסתיימים בפתרון האחרון מחלקים לזוגות (שלב i, משוב).  בסוף מאמנים מודל קטן על הדאהטסט הזה תוך משקול שונה לשלישיות והזוגות מהשלב הקודם. מטרת האימון היא חיזוי הטוקן הבא כמו שמקובל באימון מוקדם של מודלי שפה.  ככה מצליחים לשפר את הביצועים של המודל הקטן במשימות מורכבות של reasoning.' 'הם מצאו שיש שכבה מסוימת במודל להסרת הרעש שאם מעתיקים את האקטיבציות שלה עבור הקלט הטקסטואלי הלא מורעש האובייקט ״הנערך״ חוזר לתמונה המגונרטת. שימו לב שמכיוון של הארכיטקטורה של המודל מבוססת על ResNet זה התוצאה מההזנה של הקלט המורעש לא זהה לזו של לקלט הלא מורעש. אבל כן מקבלים תמונה דומה עם אותו האובייקט. וכמובן ששכבות שונות אחריות על שינוי צבע, סגנון וכדומה.' 'paragraph 13 (they found)'

# --- 9. Remove the now-duplicate empty paragraph that followed it (originally paragraph 14) ---
# After the two insertions above, the paragraph holding the text we just updated in step 8
# is at fixed index 15 (1-based); the stray empty paragraph right after it (index 16) is removed,
# leaving only the single empty paragraph before the section break.
$emptyPara = $d.Paragraphs(16)
$emptyPara.Range.Delete()

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
